$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 8: date, hours, and description ("Layout" -> new shared string)
$ws.Range("A8").Value = "1/12/2016"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "Layout"

# Update the Totaal formula in D2 to only sum the filled rows (B2:B8)
$ws.Range("D2").Formula = "=SUM(B2:B8)"

# Move the active selection to B13
$ws.Range("B13").Select() | Out-Null
